$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 408, pushing existing rows 408-433 down to 410-435.
$ws.Rows.Item(408).Insert()
$ws.Rows.Item(408).Insert()

# New row 408: Brócoli, Primera, week of 2023-04-05, Región del Maule
$ws.Cells.Item(408,1).Value = 7
$ws.Cells.Item(408,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(408,3).Value = "Ñuble"
$ws.Cells.Item(408,4).Value = 45021
$ws.Cells.Item(408,5).Value = 16
$ws.Cells.Item(408,6).Value = 100112023
$ws.Cells.Item(408,7).Value = "Brócoli"
$ws.Cells.Item(408,8).Value = "Sin especificar"
$ws.Cells.Item(408,9).Value = "Primera"
$ws.Cells.Item(408,10).Value = 100
$ws.Cells.Item(408,11).Value = 1200
$ws.Cells.Item(408,12).Value = 1200
$ws.Cells.Item(408,13).Value = 1200
$ws.Cells.Item(408,14).Value = "`$/unidad"
$ws.Cells.Item(408,15).Value = "Región del Maule"
$ws.Cells.Item(408,16).Value = 1200
$ws.Cells.Item(408,17).Value = 1
$ws.Cells.Item(408,18).Value = "Hortaliza"

# New row 409: Brócoli, Segunda, week of 2023-04-05, Región del Maule
$ws.Cells.Item(409,1).Value = 7
$ws.Cells.Item(409,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(409,3).Value = "Ñuble"
$ws.Cells.Item(409,4).Value = 45021
$ws.Cells.Item(409,5).Value = 16
$ws.Cells.Item(409,6).Value = 100112023
$ws.Cells.Item(409,7).Value = "Brócoli"
$ws.Cells.Item(409,8).Value = "Sin especificar"
$ws.Cells.Item(409,9).Value = "Segunda"
$ws.Cells.Item(409,10).Value = 100
$ws.Cells.Item(409,11).Value = 1000
$ws.Cells.Item(409,12).Value = 1000
$ws.Cells.Item(409,13).Value = 1000
$ws.Cells.Item(409,14).Value = "`$/unidad"
$ws.Cells.Item(409,15).Value = "Región del Maule"
$ws.Cells.Item(409,16).Value = 1000
$ws.Cells.Item(409,17).Value = 1
$ws.Cells.Item(409,18).Value = "Hortaliza"
